# Applies the cryptos.xlsx data refresh described in the commit diff.
# Column D values are numeric-looking text (prices like "557.22", "65.00")
# that must stay as literal text (matching the source t="inlineStr" cells),
# so we temporarily force a Text number format while assigning them, then
# restore the default "Normal" style so no stray formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "62.088.29"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +2.81%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.423.86"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +4.03%  "
$ws.Range("E4").Value = "  -0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "557.22"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.31%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "138.66"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +5.63%  "
$ws.Range("E7").Value = "  -0.05%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.583"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.16%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "2.422.35"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +4.10%  "
$ws.Range("E10").Value = "  +2.76%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.77"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +4.26%  "
$ws.Range("E12").Value = "  +0.12%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.347"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +3.61%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "25.75"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +8.88%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "2.854.54"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +3.96%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "62.023.64"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.78%  "
$ws.Range("E17").Value = "  +5.26%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.418.93"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +3.29%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.13"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +5.09%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "344.74"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +9.42%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "4.23"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.20%  "
$ws.Range("E22").Value = "  +2.84%  "
$ws.Range("E23").Value = "  +0.19%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "65.00"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("E25").Value = "  -0.06%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.23%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.52"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +12.08%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "8.29"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +5.60%  "
$ws.Range("E29").Value = "  +14.66%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.0₃0790"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +7.87%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.81"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +4.83%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "6.33"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +6.51%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "170.79"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.35%  "
$ws.Range("E34").Value = "  +4.58%  "
$ws.Range("E35").Value = "  +3.95%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "375.16"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +16.85%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "18.53"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +3.84%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "4.48"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +10.03%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("E41").Value = "  +9.46%  "
$ws.Range("E42").Value = "  +2.96%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "145.77"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +6.10%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "3.67"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +4.97%  "
$ws.Range("E45").Value = "  +8.20%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0956"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +1.84%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0519"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +4.80%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.587"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +4.26%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "18.04"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +6.47%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0221"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.44%  "
$ws.Range("E51").Value = "  +2.51%  "
